# Applies the "Updated code on 3-8-21" edit to the TestData worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Row 1: header A1 text changed from "TC" to "Tcsample"
$ws.Range("A1").Value = "Tcsample"

# Row 5: A5 now carries "TC1" (previously blank) and picks up the bold
# "header" formatting already used by its row-mates (B5:E5), matching the
# style used by the other section header rows (row 1 / row 5).
$ws.Range("A5").Value = "TC1"
$ws.Range("B5").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 7: D7 and E7 now carry "NA" (previously blank)
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "NA"

# Update the view state: scroll the window (new topLeftCell "C1") and move
# the active selection to H5 (was C11).
$ws.Activate()
$excel.Goto($ws.Range("H5"), $true)
